$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata": update Date / Contact, insert a new "Jurisdiction" row,
# and bump the concept Count from 1 to 7.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Insert a new row 11 ("Jurisdiction") - this pushes the former rows 11-21
# down to 12-22 and grows the used range to A1:B22.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$meta.Cells.Item(11, 1).Value = "Jurisdiction"
$meta.Cells.Item(11, 2).Value = ""

# Refresh the publication Date (row 8).
$meta.Cells.Item(8, 2).Value = "2024-10-02T15:04:17+00:00"

# Refresh the Contact display string (row 10).
$meta.Cells.Item(10, 2).Value = "Ferlab.bio (http://example.org/example-publisher)"

# The Count row is now row 22 (shifted down by the inserted row); update the
# concept count from 1 to 7, keeping it stored as text like the original.
$meta.Cells.Item(22, 2).NumberFormat = "@"
$meta.Cells.Item(22, 2).Value = "7"
$meta.Range("A21:B21").Copy()
$meta.Range("A22:B22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------------
# Sheet "Concepts": rename the existing "Genomics" code to lowercase
# "genomics", and append the six remaining data-category concepts.
# ---------------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

$concepts.Cells.Item(2, 2).Value = "genomics"

# Stamp rows 3-8 with the same formatting (style, borders, fill) as row 2
# before filling in their values.
$concepts.Range("A2:D2").Copy()
for ($r = 3; $r -le 8; $r++) {
    $concepts.Range("A" + $r + ":D" + $r).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$codes = @("imaging", "clinical", "transcriptomics", "proteomics", "metabolomics", "other")
$displays = @("Imaging", "Clinical", "Transcriptomics", "Proteomics", "Metabolomics", "Other")

for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 3
    $concepts.Cells.Item($r, 1).NumberFormat = "@"
    $concepts.Cells.Item($r, 1).Value = "1"
    $concepts.Cells.Item($r, 2).Value = $codes[$i]
    $concepts.Cells.Item($r, 3).Value = $displays[$i]
}

# The NumberFormat="@" stamp above creates a new style id; repaint rows 3-8
# with row 2's formatting once more so every body cell shares style "2".
$concepts.Range("A2:D2").Copy()
for ($r = 3; $r -le 8; $r++) {
    $concepts.Range("A" + $r + ":D" + $r).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}
